# Locale goldenSpiral added for all languages [en,pl,de,es,ru]
#
# The translations sheet is a simple key/value table:
#   col A = key, col B = en, col C = pl, col D = de, col E = es, col F = ru
# A new row for the "goldenSpiral" key is inserted right before the
# "selectAll" row (which currently sits at row 33), pushing every row
# below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 33; rows 33..40 (selectAll .. lineThicknessProm)
# shift down to 34..41, exactly like the OOXML diff shows.
$ws.Rows("33:33").Insert()

# Fill the freshly inserted row with the new translations.
$ws.Range("A33").Value = "goldenSpiral"
$ws.Range("B33").Value = "Golden Spiral"
$ws.Range("C33").Value = "Złota spirala"
$ws.Range("D33").Value = "Goldene Spirale"
$ws.Range("E33").Value = "Espiral Dorada"
$ws.Range("F33").Value = "Золотая спираль"

# Cosmetic touch-ups that Excel performs automatically after the edit:
# the translation columns are re-fit to the (now slightly different)
# widest content, and the view is reset to 100% zoom.
$ws.Columns("B:F").AutoFit() | Out-Null
$excel.ActiveWindow.Zoom = 100
